# Scheduled runner update: refresh market-price / leve-profit figures (columns H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 179.77777  # H9: 237.25 -> 179.77777
$ws.Cells.Item(9, 9).Value = 77.25  # I9: 83 -> 77.25
$ws.Cells.Item(9, 10).Value = 1000  # J9: 700 -> 1000
$ws.Cells.Item(9, 11).Value = 77.25  # K9: 83 -> 77.25
$ws.Cells.Item(9, 12).Value = 1000  # L9: 700 -> 1000
$ws.Cells.Item(9, 13).Value = 91.75  # M9: 86 -> 91.75
$ws.Cells.Item(9, 14).Value = -1338  # N9: -1038 -> -1338
$ws.Cells.Item(32, 8).Value = 618  # H32: 700 -> 618
$ws.Cells.Item(32, 9).Value = 382.85715  # I32: 433.33334 -> 382.85715
$ws.Cells.Item(32, 10).Value = 1166.6666  # J32: 1500 -> 1166.6666
$ws.Cells.Item(32, 11).Value = 382.85715  # K32: 433.33334 -> 382.85715
$ws.Cells.Item(32, 12).Value = 1166.6666  # L32: 1500 -> 1166.6666
$ws.Cells.Item(32, 13).Value = -56.85714999999999  # M32: -107.33334 -> -56.85714999999999
$ws.Cells.Item(32, 14).Value = -1818.6666  # N32: -2152 -> -1818.6666
$ws.Cells.Item(33, 8).Value = 291.83334  # H33: 313.9091 -> 291.83334
$ws.Cells.Item(33, 9).Value = 313.81818  # I33: 340.3 -> 313.81818
$ws.Cells.Item(33, 11).Value = 313.81818  # K33: 340.3 -> 313.81818
$ws.Cells.Item(33, 13).Value = -84.81817999999998  # M33: -111.3 -> -84.81817999999998
$ws.Cells.Item(55, 8).Value = 235.63333  # H55: 74.13793 -> 235.63333
$ws.Cells.Item(55, 9).Value = 671.125  # I55: 65 -> 671.125
$ws.Cells.Item(55, 10).Value = 77.27273  # J55: 76.521736 -> 77.27273
$ws.Cells.Item(55, 11).Value = 671.125  # K55: 65 -> 671.125
$ws.Cells.Item(55, 12).Value = 77.27273  # L55: 76.521736 -> 77.27273
$ws.Cells.Item(55, 13).Value = -457.125  # M55: 149 -> -457.125
$ws.Cells.Item(55, 14).Value = -505.27273  # N55: -504.521736 -> -505.27273
$ws.Cells.Item(112, 8).Value = 3368052.5  # H112: 3269020.8 -> 3368052.5
$ws.Cells.Item(112, 9).Value = 0  # I112: 969 -> 0
$ws.Cells.Item(112, 11).Value = 0  # K112: 2907 -> 0
$ws.Cells.Item(112, 13).ClearContents()  # M112: -1799 -> (removed)
$ws.Cells.Item(129, 8).Value = 245063.05  # H129: 223355.22 -> 245063.05
$ws.Cells.Item(129, 10).Value = 264387.75  # J129: 239288.92 -> 264387.75
$ws.Cells.Item(129, 12).Value = 793163.25  # L129: 717866.76 -> 793163.25
$ws.Cells.Item(129, 14).Value = -803163.25  # N129: -727866.76 -> -803163.25
$ws.Cells.Item(132, 8).Value = 3423.7932  # H132: 3374.2068 -> 3423.7932
$ws.Cells.Item(132, 9).Value = 3691.348  # I132: 3434.04 -> 3691.348
$ws.Cells.Item(132, 10).Value = 2398.1667  # J132: 3000.25 -> 2398.1667
$ws.Cells.Item(132, 11).Value = 11074.044  # K132: 10302.12 -> 11074.044
$ws.Cells.Item(132, 12).Value = 7194.500100000001  # L132: 9000.75 -> 7194.500100000001
$ws.Cells.Item(132, 13).Value = -8544.044  # M132: -7772.119999999999 -> -8544.044
$ws.Cells.Item(132, 14).Value = -12254.5001  # N132: -14060.75 -> -12254.5001
$ws.Cells.Item(135, 8).Value = 13161654  # H135: 20839202 -> 13161654
$ws.Cells.Item(135, 9).Value = 751.5  # I135: 1177 -> 751.5
$ws.Cells.Item(135, 10).Value = 50012184  # J135: 62515252 -> 50012184
$ws.Cells.Item(135, 11).Value = 6763.5  # K135: 10593 -> 6763.5
$ws.Cells.Item(135, 12).Value = 450109656  # L135: 562637268 -> 450109656
$ws.Cells.Item(135, 13).Value = -4228.5  # M135: -8058 -> -4228.5
$ws.Cells.Item(135, 14).Value = -450114726  # N135: -562642338 -> -450114726
$ws.Cells.Item(137, 8).Value = 1574.4147  # H137: 1989.5518 -> 1574.4147
$ws.Cells.Item(137, 9).Value = 1412.5278  # I137: 1833.2084 -> 1412.5278
$ws.Cells.Item(137, 11).Value = 4237.5834  # K137: 5499.6252 -> 4237.5834
$ws.Cells.Item(137, 13).Value = -1687.5834  # M137: -2949.6252 -> -1687.5834
$ws.Cells.Item(138, 8).Value = 10528530  # H138: 10871805 -> 10528530
$ws.Cells.Item(138, 9).Value = 21739958  # I138: 22223090 -> 21739958
$ws.Cells.Item(138, 10).Value = 3515.2856  # J138: 3552.5532 -> 3515.2856
$ws.Cells.Item(138, 11).Value = 65219874  # K138: 66669270 -> 65219874
$ws.Cells.Item(138, 12).Value = 10545.8568  # L138: 10657.6596 -> 10545.8568
$ws.Cells.Item(138, 13).Value = -65214734  # M138: -66664130 -> -65214734
$ws.Cells.Item(138, 14).Value = -20825.8568  # N138: -20937.6596 -> -20825.8568
$ws.Cells.Item(141, 8).Value = 1035.1587  # H141: 1249.92 -> 1035.1587
$ws.Cells.Item(141, 9).Value = 725.11536  # I141: 837.02325 -> 725.11536
$ws.Cells.Item(141, 10).Value = 2500.818  # J141: 3786.2856 -> 2500.818
$ws.Cells.Item(141, 11).Value = 2175.34608  # K141: 2511.06975 -> 2175.34608
$ws.Cells.Item(141, 12).Value = 7502.454000000001  # L141: 11358.8568 -> 7502.454000000001
$ws.Cells.Item(141, 13).Value = 3004.65392  # M141: 2668.93025 -> 3004.65392
$ws.Cells.Item(141, 14).Value = -17862.454  # N141: -21718.8568 -> -17862.454

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4015.9658  # H32: 3926 -> 4015.9658
$ws.Cells.Item(32, 9).Value = 2654.9  # I32: 2589.3538 -> 2654.9
$ws.Cells.Item(32, 11).Value = 2654.9  # K32: 2589.3538 -> 2654.9
$ws.Cells.Item(32, 13).Value = -2367.9  # M32: -2302.3538 -> -2367.9
$ws.Cells.Item(45, 8).Value = 3422.3794  # H45: 4956.909 -> 3422.3794
$ws.Cells.Item(45, 9).Value = 3053.1  # I45: 5385.3335 -> 3053.1
$ws.Cells.Item(45, 10).Value = 4243  # J45: 4442.8 -> 4243
$ws.Cells.Item(45, 11).Value = 3053.1  # K45: 5385.3335 -> 3053.1
$ws.Cells.Item(45, 12).Value = 4243  # L45: 4442.8 -> 4243
$ws.Cells.Item(45, 13).Value = -2676.1  # M45: -5008.3335 -> -2676.1
$ws.Cells.Item(45, 14).Value = -4997  # N45: -5196.8 -> -4997
$ws.Cells.Item(61, 8).Value = 581911.8  # H61: 347282.53 -> 581911.8
$ws.Cells.Item(61, 9).Value = 601275.1  # I61: 429563.72 -> 601275.1
$ws.Cells.Item(61, 10).Value = 1014  # J61: 1701.5 -> 1014
$ws.Cells.Item(61, 11).Value = 601275.1  # K61: 429563.72 -> 601275.1
$ws.Cells.Item(61, 12).Value = 1014  # L61: 1701.5 -> 1014
$ws.Cells.Item(61, 13).Value = -601063.1  # M61: -429351.72 -> -601063.1
$ws.Cells.Item(61, 14).Value = -1438  # N61: -2125.5 -> -1438
$ws.Cells.Item(74, 8).Value = 24392068  # H74: 25642922 -> 24392068
$ws.Cells.Item(74, 9).Value = 26317524  # I74: 27779586 -> 26317524
$ws.Cells.Item(74, 11).Value = 26317524  # K74: 27779586 -> 26317524
$ws.Cells.Item(74, 13).Value = -26316650  # M74: -27778712 -> -26316650
$ws.Cells.Item(77, 8).Value = 24392068  # H77: 25642922 -> 24392068
$ws.Cells.Item(77, 9).Value = 26317524  # I77: 27779586 -> 26317524
$ws.Cells.Item(77, 11).Value = 131587620  # K77: 138897930 -> 131587620
$ws.Cells.Item(77, 13).Value = -131583252  # M77: -138893562 -> -131583252
$ws.Cells.Item(132, 8).Value = 14463.154  # H132: 9346.229499999999 -> 14463.154
$ws.Cells.Item(132, 9).Value = 1610.5758  # I132: 1096.7778 -> 1610.5758
$ws.Cells.Item(132, 10).Value = 85152.336  # J132: 72984.86 -> 85152.336
$ws.Cells.Item(132, 11).Value = 4831.7274  # K132: 3290.3334 -> 4831.7274
$ws.Cells.Item(132, 12).Value = 255457.008  # L132: 218954.58 -> 255457.008
$ws.Cells.Item(132, 13).Value = -2301.7274  # M132: -760.3334000000004 -> -2301.7274
$ws.Cells.Item(132, 14).Value = -260517.008  # N132: -224014.58 -> -260517.008
$ws.Cells.Item(136, 8).Value = 581911.8  # H136: 347282.53 -> 581911.8
$ws.Cells.Item(136, 9).Value = 601275.1  # I136: 429563.72 -> 601275.1
$ws.Cells.Item(136, 10).Value = 1014  # J136: 1701.5 -> 1014
$ws.Cells.Item(136, 11).Value = 1803825.3  # K136: 1288691.16 -> 1803825.3
$ws.Cells.Item(136, 12).Value = 3042  # L136: 5104.5 -> 3042
$ws.Cells.Item(136, 13).Value = -1801275.3  # M136: -1286141.16 -> -1801275.3
$ws.Cells.Item(136, 14).Value = -8142  # N136: -10204.5 -> -8142

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3036.383  # H134: 3364.976 -> 3036.383
$ws.Cells.Item(134, 9).Value = 3473.0571  # I134: 3772.9688 -> 3473.0571
$ws.Cells.Item(134, 10).Value = 1762.75  # J134: 2059.4 -> 1762.75
$ws.Cells.Item(134, 11).Value = 10419.1713  # K134: 11318.9064 -> 10419.1713
$ws.Cells.Item(134, 12).Value = 5288.25  # L134: 6178.200000000001 -> 5288.25
$ws.Cells.Item(134, 13).Value = -7884.1713  # M134: -8783.9064 -> -7884.1713
$ws.Cells.Item(134, 14).Value = -10358.25  # N134: -11248.2 -> -10358.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 16602.908  # H58: 8637.216 -> 16602.908
$ws.Cells.Item(58, 9).Value = 1255.1305  # I58: 808.1395 -> 1255.1305
$ws.Cells.Item(58, 10).Value = 51902.8  # J58: 23939.5 -> 51902.8
$ws.Cells.Item(58, 11).Value = 1255.1305  # K58: 808.1395 -> 1255.1305
$ws.Cells.Item(58, 12).Value = 51902.8  # L58: 23939.5 -> 51902.8
$ws.Cells.Item(58, 13).Value = -1052.1305  # M58: -605.1395 -> -1052.1305
$ws.Cells.Item(58, 14).Value = -52308.8  # N58: -24345.5 -> -52308.8
$ws.Cells.Item(68, 8).Value = 51206.332  # H68: 47154.75 -> 51206.332
$ws.Cells.Item(68, 10).Value = 51206.332  # J68: 47154.75 -> 51206.332
$ws.Cells.Item(68, 12).Value = 51206.332  # L68: 47154.75 -> 51206.332
$ws.Cells.Item(68, 14).Value = -52704.332  # N68: -48652.75 -> -52704.332
$ws.Cells.Item(71, 8).Value = 51206.332  # H71: 47154.75 -> 51206.332
$ws.Cells.Item(71, 10).Value = 51206.332  # J71: 47154.75 -> 51206.332
$ws.Cells.Item(71, 12).Value = 153618.996  # L71: 141464.25 -> 153618.996
$ws.Cells.Item(71, 14).Value = -161106.996  # N71: -148952.25 -> -161106.996
$ws.Cells.Item(132, 8).Value = 1779.8541  # H132: 2409.625 -> 1779.8541
$ws.Cells.Item(132, 9).Value = 1383.0217  # I132: 1809.4839 -> 1383.0217
$ws.Cells.Item(132, 10).Value = 10907  # J132: 21014 -> 10907
$ws.Cells.Item(132, 11).Value = 4149.0651  # K132: 5428.4517 -> 4149.0651
$ws.Cells.Item(132, 12).Value = 32721  # L132: 63042 -> 32721
$ws.Cells.Item(132, 13).Value = -1619.0651  # M132: -2898.4517 -> -1619.0651
$ws.Cells.Item(132, 14).Value = -37781  # N132: -68102 -> -37781
$ws.Cells.Item(134, 8).Value = 1066.8788  # H134: 1193.6522 -> 1066.8788
$ws.Cells.Item(134, 9).Value = 965.1539  # I134: 1023.9474 -> 965.1539
$ws.Cells.Item(134, 10).Value = 1444.7142  # J134: 1999.75 -> 1444.7142
$ws.Cells.Item(134, 11).Value = 2895.4617  # K134: 3071.8422 -> 2895.4617
$ws.Cells.Item(134, 12).Value = 4334.142599999999  # L134: 5999.25 -> 4334.142599999999
$ws.Cells.Item(134, 13).Value = -360.4616999999998  # M134: -536.8422 -> -360.4616999999998
$ws.Cells.Item(134, 14).Value = -9404.142599999999  # N134: -11069.25 -> -9404.142599999999
$ws.Cells.Item(136, 8).Value = 16602.908  # H136: 8637.216 -> 16602.908
$ws.Cells.Item(136, 9).Value = 1255.1305  # I136: 808.1395 -> 1255.1305
$ws.Cells.Item(136, 10).Value = 51902.8  # J136: 23939.5 -> 51902.8
$ws.Cells.Item(136, 11).Value = 3765.3915  # K136: 2424.4185 -> 3765.3915
$ws.Cells.Item(136, 12).Value = 155708.4  # L136: 71818.5 -> 155708.4
$ws.Cells.Item(136, 13).Value = -1215.3915  # M136: 125.5815000000002 -> -1215.3915
$ws.Cells.Item(136, 14).Value = -160808.4  # N136: -76918.5 -> -160808.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1480.8096  # H5: 1642.7368 -> 1480.8096
$ws.Cells.Item(5, 9).Value = 1131.5555  # I5: 1470.5714 -> 1131.5555
$ws.Cells.Item(5, 10).Value = 1742.75  # J5: 1743.1666 -> 1742.75
$ws.Cells.Item(5, 11).Value = 3394.6665  # K5: 4411.7142 -> 3394.6665
$ws.Cells.Item(5, 12).Value = 5228.25  # L5: 5229.4998 -> 5228.25
$ws.Cells.Item(5, 13).Value = -3282.6665  # M5: -4299.7142 -> -3282.6665
$ws.Cells.Item(5, 14).Value = -5452.25  # N5: -5453.4998 -> -5452.25
$ws.Cells.Item(16, 8).Value = 100  # H16: 0 -> 100
$ws.Cells.Item(16, 10).Value = 100  # J16: 0 -> 100
$ws.Cells.Item(16, 12).Value = 300  # L16: 0 -> 300
$ws.Cells.Item(16, 14).Value = -646  # N16: None -> -646
$ws.Cells.Item(39, 8).Value = 2587.3  # H39: 2793.3333 -> 2587.3
$ws.Cells.Item(39, 10).Value = 2587.3  # J39: 2793.3333 -> 2587.3
$ws.Cells.Item(39, 12).Value = 7761.900000000001  # L39: 8379.999899999999 -> 7761.900000000001
$ws.Cells.Item(39, 14).Value = -8349.900000000001  # N39: -8967.999899999999 -> -8349.900000000001
$ws.Cells.Item(55, 8).Value = 4000  # H55: 2642.8572 -> 4000
$ws.Cells.Item(55, 10).Value = 4000  # J55: 2642.8572 -> 4000
$ws.Cells.Item(55, 12).Value = 12000  # L55: 7928.571599999999 -> 12000
$ws.Cells.Item(55, 14).Value = -12354  # N55: -8282.571599999999 -> -12354
$ws.Cells.Item(122, 8).Value = 727.2727  # H122: 740 -> 727.2727
$ws.Cells.Item(122, 10).Value = 1220  # J122: 1163.3334 -> 1220
$ws.Cells.Item(122, 12).Value = 10980  # L122: 10470.0006 -> 10980
$ws.Cells.Item(122, 14).Value = -15880  # N122: -15370.0006 -> -15880
$ws.Cells.Item(131, 8).Value = 798.08  # H131: 758.48486 -> 798.08
$ws.Cells.Item(131, 10).Value = 811.76843  # J131: 770.21277 -> 811.76843
$ws.Cells.Item(131, 12).Value = 2435.30529  # L131: 2310.63831 -> 2435.30529
$ws.Cells.Item(131, 14).Value = -12515.30529  # N131: -12390.63831 -> -12515.30529
$ws.Cells.Item(135, 8).Value = 1480.8096  # H135: 1642.7368 -> 1480.8096
$ws.Cells.Item(135, 9).Value = 1131.5555  # I135: 1470.5714 -> 1131.5555
$ws.Cells.Item(135, 10).Value = 1742.75  # J135: 1743.1666 -> 1742.75
$ws.Cells.Item(135, 11).Value = 10183.9995  # K135: 13235.1426 -> 10183.9995
$ws.Cells.Item(135, 12).Value = 15684.75  # L135: 15688.4994 -> 15684.75
$ws.Cells.Item(135, 13).Value = -7648.9995  # M135: -10700.1426 -> -7648.9995
$ws.Cells.Item(135, 14).Value = -20754.75  # N135: -20758.4994 -> -20754.75
$ws.Cells.Item(141, 8).Value = 3775  # H141: 4049.5454 -> 3775
$ws.Cells.Item(141, 9).Value = 3549.0908  # I141: 3834.4443 -> 3549.0908
$ws.Cells.Item(141, 11).Value = 10647.2724  # K141: 11503.3329 -> 10647.2724
$ws.Cells.Item(141, 13).Value = -5467.2724  # M141: -6323.332900000001 -> -5467.2724

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(111, 8).Value = 29999  # H111: 30000 -> 29999
$ws.Cells.Item(111, 10).Value = 29999  # J111: 30000 -> 29999
$ws.Cells.Item(111, 12).Value = 29999  # L111: 30000 -> 29999
$ws.Cells.Item(111, 14).Value = -36133  # N111: -36134 -> -36133

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1048.0385  # H132: 1056.3529 -> 1048.0385
$ws.Cells.Item(132, 9).Value = 1048.0385  # I132: 1067.42 -> 1048.0385
$ws.Cells.Item(132, 10).Value = 0  # J132: 503 -> 0
$ws.Cells.Item(132, 11).Value = 3144.1155  # K132: 3202.26 -> 3144.1155
$ws.Cells.Item(132, 12).Value = 0  # L132: 1509 -> 0
$ws.Cells.Item(132, 13).Value = -614.1155000000003  # M132: -672.2600000000002 -> -614.1155000000003
$ws.Cells.Item(132, 14).ClearContents()  # N132: -6569 -> (removed)
$ws.Cells.Item(136, 8).Value = 946.6429000000001  # H136: 828.902 -> 946.6429000000001
$ws.Cells.Item(136, 9).Value = 896.17645  # I136: 767.8409 -> 896.17645
$ws.Cells.Item(136, 10).Value = 1161.125  # J136: 1212.7142 -> 1161.125
$ws.Cells.Item(136, 11).Value = 2688.52935  # K136: 2303.5227 -> 2688.52935
$ws.Cells.Item(136, 12).Value = 3483.375  # L136: 3638.1426 -> 3483.375
$ws.Cells.Item(136, 13).Value = -138.5293500000002  # M136: 246.4773 -> -138.5293500000002
$ws.Cells.Item(136, 14).Value = -8583.375  # N136: -8738.142599999999 -> -8583.375

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 560.0244  # H132: 700.55817 -> 560.0244
$ws.Cells.Item(132, 9).Value = 572.525  # I132: 636.3611 -> 572.525
$ws.Cells.Item(132, 10).Value = 60  # J132: 1030.7142 -> 60
$ws.Cells.Item(132, 11).Value = 1717.575  # K132: 1909.0833 -> 1717.575
$ws.Cells.Item(132, 12).Value = 180  # L132: 3092.1426 -> 180
$ws.Cells.Item(132, 13).Value = 812.4250000000002  # M132: 620.9167000000002 -> 812.4250000000002
$ws.Cells.Item(132, 14).Value = -5240  # N132: -8152.142599999999 -> -5240
$ws.Cells.Item(136, 8).Value = 16668693  # H136: 18183952 -> 16668693
$ws.Cells.Item(136, 9).Value = 23256748  # I136: 25642026 -> 23256748
$ws.Cells.Item(136, 10).Value = 4788.7646  # J136: 4900.5625 -> 4788.7646
$ws.Cells.Item(136, 11).Value = 69770244  # K136: 76926078 -> 69770244
$ws.Cells.Item(136, 12).Value = 14366.2938  # L136: 14701.6875 -> 14366.2938
$ws.Cells.Item(136, 13).Value = -69767694  # M136: -76923528 -> -69767694
$ws.Cells.Item(136, 14).Value = -19466.2938  # N136: -19801.6875 -> -19466.2938

